$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Clone the formatting (bold headers + borders) of the existing
#        A1:E4 result table onto the new G1:K4 block --------------------
$ws.Range("A1:E4").Copy()
$ws.Range("G1").PasteSpecial(-4122)   # xlPasteFormats

# --- Helper: write an arbitrary (possibly numeric-looking) string into a
#     cell as genuine TEXT, the same way the existing "6858.28"-style
#     score labels are stored, without leaving any trace behind.
#     (A scratch cell is formatted as Text, filled, copied, and its VALUE
#     ONLY is pasted onto the destination so destination formatting from
#     step 1 above is preserved; the scratch cell is fully cleared after.)
function Set-TextValue($cellRef, $text) {
    $ws.Range("ZZ1").NumberFormat = "@"
    $ws.Range("ZZ1").Value = $text
    $ws.Range("ZZ1").Copy()
    $ws.Range($cellRef).PasteSpecial(-4163)   # xlPasteValues
    $ws.Range("ZZ1").Clear()
}

# --- 2. Header row (row 1): label + the 4 new budget columns -----------
Set-TextValue "G1" "max score"
$ws.Range("H1").Value = 600
$ws.Range("I1").Value = 700
$ws.Range("J1").Value = 800
$ws.Range("K1").Value = 900

# --- 3. Row labels (time limits) in column G, stored as text like the
#        original 0.4 / 0.45 / 0.5 / 0.55 labels -------------------------
Set-TextValue "G2" "0.4125"
Set-TextValue "G3" "0.425"
Set-TextValue "G4" "0.4375"

# --- 4. New experiment results (stored as text, matching the existing
#        score cells) ----------------------------------------------------
Set-TextValue "H2" "6810.28"
Set-TextValue "I2" "6799.84"
Set-TextValue "J2" "6835.92"
Set-TextValue "K2" "6830.92"

Set-TextValue "H3" "6804.92"
Set-TextValue "I3" "6802.28"
Set-TextValue "J3" "6816.20"
Set-TextValue "K3" "6817.92"

Set-TextValue "H4" "6762.28"
Set-TextValue "I4" "6792.92"
Set-TextValue "J4" "6822.28"
Set-TextValue "K4" "6785.56"

# --- 5. Match the author's final selection -------------------------------
$ws.Range("K4").Select()
